$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$ws.Range("B9").Value = 35697.057505597
$ws.Range("D9").Value = 151.115833761224
$ws.Range("F9").Value = 15176.937152828
$ws.Range("H9").Value = 74.2330705107738
$ws.Range("J9").Value = 95953.019481924
$ws.Range("L9").Value = 168780.715838675
$ws.Range("M9").Value = 5.91080491354399
$ws.Range("O9").Value = 10.8426965476467
